# [Fonds de solidarite] Add 2020-08-31 data
#
# The sheet holds "nombre_aides" (col C) and "montant_total" (col D) as
# text values (Excel shows no green-triangle "number stored as text"
# conversion is wanted here — the source feed keeps them as strings), so
# each new value is entered with a leading apostrophe to force text entry
# and keep it from being re-typed as a Number by Excel's automatic
# value-coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corse (reg 94)
$ws.Range("C27").Value = "'46"
$ws.Range("D27").Value = "'123419.13"

$ws.Range("C28").Value = "'138"
$ws.Range("D28").Value = "'426965.56"

$ws.Range("C29").Value = "'60"
$ws.Range("D29").Value = "'230297.00"

$ws.Range("C31").Value = "'4"
$ws.Range("D31").Value = "'22000.00"

$ws.Range("C32").Value = "'17"
$ws.Range("D32").Value = "'34000.00"

# Grand Est (reg 44)
$ws.Range("C33").Value = "'110"
$ws.Range("D33").Value = "'319173.00"

$ws.Range("C34").Value = "'578"
$ws.Range("D34").Value = "'1907271.66"

$ws.Range("C35").Value = "'238"
$ws.Range("D35").Value = "'1225868.11"

# Guyane (reg 03)
$ws.Range("C45").Value = "'27"
$ws.Range("D45").Value = "'102621.84"

$ws.Range("C46").Value = "'91"
$ws.Range("D46").Value = "'409774.61"

$ws.Range("C47").Value = "'50"
$ws.Range("D47").Value = "'299203.00"

$ws.Range("C49").Value = "'5"
$ws.Range("D49").Value = "'26000.00"

# Hauts-de-France (reg 32)
$ws.Range("C51").Value = "'106"
$ws.Range("D51").Value = "'317768.17"

$ws.Range("C52").Value = "'609"
$ws.Range("D52").Value = "'2207241.21"

$ws.Range("C53").Value = "'267"
$ws.Range("D53").Value = "'1192878.76"

$ws.Range("C54").Value = "'92"
$ws.Range("D54").Value = "'548274.23"

$ws.Range("C55").Value = "'27"
$ws.Range("D55").Value = "'153213.00"

# Normandie (reg 28)
$ws.Range("C76").Value = "'98"
$ws.Range("D76").Value = "'266139.87"

$ws.Range("C77").Value = "'410"
$ws.Range("D77").Value = "'1338019.84"

$ws.Range("C78").Value = "'162"
$ws.Range("D78").Value = "'672987.18"

$ws.Range("C79").Value = "'46"
$ws.Range("D79").Value = "'222621.67"

$ws.Range("C80").Value = "'12"
$ws.Range("D80").Value = "'81000.00"

$ws.Range("C81").Value = "'16"
$ws.Range("D81").Value = "'32000.00"
